$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 304.4375
$ws.Range("I6").Value = 307.69232
$ws.Range("K6").Value = 923.07696
$ws.Range("M6").Value = -811.07696
$ws.Range("H17").Value = 695.7213
$ws.Range("J17").Value = 695.7213
$ws.Range("L17").Value = 2087.1639
$ws.Range("N17").Value = -2423.1639
$ws.Range("H33").Value = 368.4
$ws.Range("I33").Value = 151.125
$ws.Range("K33").Value = 151.125
$ws.Range("M33").Value = 77.875
$ws.Range("H41").Value = 366
$ws.Range("I41").Value = 567
$ws.Range("J41").Value = 165
$ws.Range("K41").Value = 567
$ws.Range("L41").Value = 165
$ws.Range("M41").Value = -127
$ws.Range("N41").Value = -1045
$ws.Range("H63").Value = 75000
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 75000
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H69").Value = 10775
$ws.Range("I69").Value = 7413.125
$ws.Range("J69").Value = 17498.75
$ws.Range("K69").Value = 22239.375
$ws.Range("L69").Value = 52496.25
$ws.Range("M69").Value = -21365.375
$ws.Range("N69").Value = -54244.25
$ws.Range("H72").Value = 10775
$ws.Range("I72").Value = 7413.125
$ws.Range("J72").Value = 17498.75
$ws.Range("K72").Value = 66718.125
$ws.Range("L72").Value = 157488.75
$ws.Range("M72").Value = -62350.125
$ws.Range("N72").Value = -166224.75
$ws.Range("H86").Value = 3796.4736
$ws.Range("I86").Value = 3699.2856
$ws.Range("K86").Value = 3699.2856
$ws.Range("M86").Value = -2576.2856
$ws.Range("H89").Value = 3796.4736
$ws.Range("I89").Value = 3699.2856
$ws.Range("K89").Value = 18496.428
$ws.Range("M89").Value = -12880.428
$ws.Range("H108").Value = 60000
$ws.Range("J108").Value = 60000
$ws.Range("L108").Value = 60000
$ws.Range("N108").Value = -67680
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H115").Value = 832.6667
$ws.Range("J115").Value = 500
$ws.Range("L115").Value = 1500
$ws.Range("N115").Value = -4634
$ws.Range("H116").Value = 15731.346
$ws.Range("I116").Value = 15078.375
$ws.Range("J116").Value = 16021.556
$ws.Range("K116").Value = 15078.375
$ws.Range("L116").Value = 16021.556
$ws.Range("M116").Value = -11636.375
$ws.Range("N116").Value = -22905.556
$ws.Range("H137").Value = 2766.7812
$ws.Range("I137").Value = 2289.1052
$ws.Range("J137").Value = 3464.923
$ws.Range("K137").Value = 6867.3156
$ws.Range("L137").Value = 10394.769
$ws.Range("M137").Value = -4317.3156
$ws.Range("N137").Value = -15494.769
$ws.Range("H138").Value = 20325.182
$ws.Range("I138").Value = 25079.424
$ws.Range("J138").Value = 2666.5715
$ws.Range("K138").Value = 75238.272
$ws.Range("L138").Value = 7999.7145
$ws.Range("M138").Value = -70098.272
$ws.Range("N138").Value = -18279.7145

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12155.06
$ws.Range("I32").Value = 12307.656
$ws.Range("K32").Value = 12307.656
$ws.Range("M32").Value = -12020.656
$ws.Range("H45").Value = 3535.5
$ws.Range("I45").Value = 2161.3333
$ws.Range("K45").Value = 2161.3333
$ws.Range("M45").Value = -1784.3333
$ws.Range("H61").Value = 11950.066
$ws.Range("I61").Value = 15454.5
$ws.Range("K61").Value = 15454.5
$ws.Range("M61").Value = -15242.5
$ws.Range("H74").Value = 899.8
$ws.Range("I74").Value = 859.76
$ws.Range("J74").Value = 1100
$ws.Range("K74").Value = 859.76
$ws.Range("L74").Value = 1100
$ws.Range("M74").Value = 14.24000000000001
$ws.Range("N74").Value = -2848
$ws.Range("H77").Value = 899.8
$ws.Range("I77").Value = 859.76
$ws.Range("J77").Value = 1100
$ws.Range("K77").Value = 4298.8
$ws.Range("L77").Value = 5500
$ws.Range("M77").Value = 69.19999999999982
$ws.Range("N77").Value = -14236
$ws.Range("H97").Value = 2303.2
$ws.Range("J97").Value = 3991.125
$ws.Range("L97").Value = 3991.125
$ws.Range("N97").Value = -4983.125
$ws.Range("H102").Value = 6487.3335
$ws.Range("I102").Value = 6923.5
$ws.Range("K102").Value = 6923.5
$ws.Range("M102").Value = -5301.5
$ws.Range("H122").Value = 2337.842
$ws.Range("J122").Value = 2753.5
$ws.Range("L122").Value = 8260.5
$ws.Range("N122").Value = -13160.5
$ws.Range("H131").Value = 79995.5
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 79995.5
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 79995.5
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -90075.5
$ws.Range("H132").Value = 37602.31
$ws.Range("I132").Value = 52701.05
$ws.Range("J132").Value = 4049.5557
$ws.Range("K132").Value = 158103.15
$ws.Range("L132").Value = 12148.6671
$ws.Range("M132").Value = -155573.15
$ws.Range("N132").Value = -17208.6671
$ws.Range("H136").Value = 11950.066
$ws.Range("I136").Value = 15454.5
$ws.Range("K136").Value = 46363.5
$ws.Range("M136").Value = -43813.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H86").Value = 3399.2666
$ws.Range("I86").Value = 2379.6
$ws.Range("J86").Value = 4418.933
$ws.Range("K86").Value = 2379.6
$ws.Range("L86").Value = 4418.933
$ws.Range("M86").Value = -1256.6
$ws.Range("N86").Value = -6664.933
$ws.Range("H89").Value = 3399.2666
$ws.Range("I89").Value = 2379.6
$ws.Range("J89").Value = 4418.933
$ws.Range("K89").Value = 11898
$ws.Range("L89").Value = 22094.665
$ws.Range("M89").Value = -6282
$ws.Range("N89").Value = -33326.665
$ws.Range("H99").Value = 3882.4119
$ws.Range("I99").Value = 2849.9
$ws.Range("J99").Value = 5357.4287
$ws.Range("K99").Value = 2849.9
$ws.Range("L99").Value = 5357.4287
$ws.Range("M99").Value = -1351.9
$ws.Range("N99").Value = -8353.4287
$ws.Range("H105").Value = 3540.5334
$ws.Range("I105").Value = 3377.16
$ws.Range("J105").Value = 4357.4
$ws.Range("K105").Value = 3377.16
$ws.Range("L105").Value = 4357.4
$ws.Range("M105").Value = -1630.16
$ws.Range("N105").Value = -7851.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H31").Value = 1425.4762
$ws.Range("I31").Value = 1444.75
$ws.Range("J31").Value = 1363.8
$ws.Range("K31").Value = 1444.75
$ws.Range("L31").Value = 1363.8
$ws.Range("M31").Value = -1149.75
$ws.Range("N31").Value = -1953.8
$ws.Range("H34").Value = 1425.4762
$ws.Range("I34").Value = 1444.75
$ws.Range("J34").Value = 1363.8
$ws.Range("K34").Value = 1444.75
$ws.Range("L34").Value = 1363.8
$ws.Range("M34").Value = -1242.75
$ws.Range("N34").Value = -1767.8
$ws.Range("H86").Value = 9998.5
$ws.Range("I86").Value = 9998.5
$ws.Range("K86").Value = 9998.5
$ws.Range("M86").Value = -8875.5
$ws.Range("H89").Value = 9998.5
$ws.Range("I89").Value = 9998.5
$ws.Range("K89").Value = 49992.5
$ws.Range("M89").Value = -44376.5
$ws.Range("H93").Value = 21220
$ws.Range("I93").Value = 18022.223
$ws.Range("J93").Value = 50000
$ws.Range("K93").Value = 18022.223
$ws.Range("L93").Value = 50000
$ws.Range("M93").Value = -16150.223
$ws.Range("N93").Value = -53744
$ws.Range("H132").Value = 2594.6191
$ws.Range("I132").Value = 2472.2856
$ws.Range("J132").Value = 2655.7856
$ws.Range("K132").Value = 7416.8568
$ws.Range("L132").Value = 7967.3568
$ws.Range("M132").Value = -4886.8568
$ws.Range("N132").Value = -13027.3568
$ws.Range("H134").Value = 57690.723
$ws.Range("I134").Value = 92685.09
$ws.Range("J134").Value = 2699.5715
$ws.Range("K134").Value = 278055.27
$ws.Range("L134").Value = 8098.7145
$ws.Range("M134").Value = -275520.27
$ws.Range("N134").Value = -13168.7145
$ws.Range("H137").Value = 60997.6
$ws.Range("J137").Value = 99994
$ws.Range("L137").Value = 99994
$ws.Range("N137").Value = -110194

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 77.181816
$ws.Range("I2").Value = 50
$ws.Range("J2").Value = 149.66667
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 898.0000200000001
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -1124.00002
$ws.Range("H7").Value = 2092.6667
$ws.Range("I7").Value = 317.5
$ws.Range("J7").Value = 3512.8
$ws.Range("K7").Value = 952.5
$ws.Range("L7").Value = 10538.4
$ws.Range("M7").Value = -840.5
$ws.Range("N7").Value = -10762.4
$ws.Range("H8").Value = 640.0625
$ws.Range("I8").Value = 640.0625
$ws.Range("K8").Value = 1920.1875
$ws.Range("M8").Value = -1781.1875
$ws.Range("H32").Value = 4500
$ws.Range("J32").Value = 4500
$ws.Range("L32").Value = 13500
$ws.Range("N32").Value = -14066
$ws.Range("H38").Value = 501.03333
$ws.Range("J38").Value = 462.2857
$ws.Range("L38").Value = 1386.8571
$ws.Range("N38").Value = -2080.8571
$ws.Range("H75").Value = 198.33333
$ws.Range("J75").Value = 196
$ws.Range("L75").Value = 588
$ws.Range("N75").Value = -2584
$ws.Range("H78").Value = 198.33333
$ws.Range("J78").Value = 196
$ws.Range("L78").Value = 1764
$ws.Range("N78").Value = -11748
$ws.Range("H81").Value = 7800
$ws.Range("I81").Value = 600
$ws.Range("K81").Value = 1800
$ws.Range("M81").Value = -677
$ws.Range("H84").Value = 7800
$ws.Range("I84").Value = 600
$ws.Range("K84").Value = 5400
$ws.Range("M84").Value = 216
$ws.Range("H86").Value = 579.25
$ws.Range("I86").Value = 472.33334
$ws.Range("K86").Value = 1417.00002
$ws.Range("M86").Value = -231.0000199999999
$ws.Range("H89").Value = 579.25
$ws.Range("I89").Value = 472.33334
$ws.Range("K89").Value = 4251.00006
$ws.Range("M89").Value = 1676.99994
$ws.Range("H98").Value = 3999.5
$ws.Range("J98").Value = 3999.5
$ws.Range("L98").Value = 11998.5
$ws.Range("N98").Value = -14994.5
$ws.Range("H113").Value = 547.8570999999999
$ws.Range("I113").Value = 624.3333
$ws.Range("K113").Value = 1872.9999
$ws.Range("M113").Value = 297.0001
$ws.Range("H121").Value = 30
$ws.Range("I121").Value = 30
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 90
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = 1220
$ws.Range("N121").ClearContents()
$ws.Range("H122").Value = 710.2308
$ws.Range("I122").Value = 462.33334
$ws.Range("K122").Value = 4161.00006
$ws.Range("M122").Value = -1711.00006
$ws.Range("H131").Value = 23224.555
$ws.Range("I131").Value = 1005.25
$ws.Range("J131").Value = 41000
$ws.Range("K131").Value = 3015.75
$ws.Range("L131").Value = 123000
$ws.Range("M131").Value = 2024.25
$ws.Range("N131").Value = -133080
$ws.Range("H136").Value = 5798.8
$ws.Range("I136").Value = 5798.8
$ws.Range("K136").Value = 17396.4
$ws.Range("M136").Value = -12296.4
$ws.Range("H140").Value = 4818.0835
$ws.Range("I140").Value = 4818.0835
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 14454.2505
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -9274.250499999998
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 300
$ws.Range("I141").Value = 300
$ws.Range("K141").Value = 900
$ws.Range("M141").Value = 4280

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H80").Value = 3857.5715
$ws.Range("I80").Value = 3250.75
$ws.Range("J80").Value = 4666.6665
$ws.Range("K80").Value = 3250.75
$ws.Range("L80").Value = 4666.6665
$ws.Range("M80").Value = -2252.75
$ws.Range("N80").Value = -6662.6665
$ws.Range("H83").Value = 3857.5715
$ws.Range("I83").Value = 3250.75
$ws.Range("J83").Value = 4666.6665
$ws.Range("K83").Value = 16253.75
$ws.Range("L83").Value = 23333.3325
$ws.Range("M83").Value = -11261.75
$ws.Range("N83").Value = -33317.3325
$ws.Range("H97").Value = 1174.875
$ws.Range("I97").Value = 922.0909
$ws.Range("J97").Value = 3955.5
$ws.Range("K97").Value = 922.0909
$ws.Range("L97").Value = 3955.5
$ws.Range("M97").Value = -426.0909
$ws.Range("N97").Value = -4947.5
$ws.Range("H102").Value = 4508.353
$ws.Range("I102").Value = 3531.6365
$ws.Range("J102").Value = 6299
$ws.Range("K102").Value = 3531.6365
$ws.Range("L102").Value = 6299
$ws.Range("M102").Value = -1909.6365
$ws.Range("N102").Value = -9543
$ws.Range("H122").Value = 3770.238
$ws.Range("I122").Value = 2365.2727
$ws.Range("K122").Value = 7095.8181
$ws.Range("M122").Value = -4645.8181
$ws.Range("H126").Value = 5260.968
$ws.Range("I126").Value = 4574.478
$ws.Range("K126").Value = 13723.434
$ws.Range("M126").Value = -11253.434
$ws.Range("H132").Value = 203169
$ws.Range("I132").Value = 203169
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 609507
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -606977
$ws.Range("N132").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2080.625
$ws.Range("I16").Value = 2078.0356
$ws.Range("K16").Value = 2078.0356
$ws.Range("M16").Value = -1908.0356
$ws.Range("H64").Value = 21197
$ws.Range("J64").Value = 21197
$ws.Range("L64").Value = 21197
$ws.Range("N64").Value = -21647
$ws.Range("H67").Value = 21197
$ws.Range("J67").Value = 21197
$ws.Range("L67").Value = 21197
$ws.Range("N67").Value = -22757
$ws.Range("H68").Value = 5839.4
$ws.Range("I68").Value = 4250.3335
$ws.Range("J68").Value = 8223
$ws.Range("K68").Value = 4250.3335
$ws.Range("L68").Value = 8223
$ws.Range("M68").Value = -3501.3335
$ws.Range("N68").Value = -9721
$ws.Range("H71").Value = 5839.4
$ws.Range("I71").Value = 4250.3335
$ws.Range("J71").Value = 8223
$ws.Range("K71").Value = 21251.6675
$ws.Range("L71").Value = 41115
$ws.Range("M71").Value = -17507.6675
$ws.Range("N71").Value = -48603
$ws.Range("H88").Value = 850031.4399999999
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 850031.4399999999
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 850031.4399999999
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -850887.4399999999
$ws.Range("H91").Value = 850031.4399999999
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 850031.4399999999
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 850031.4399999999
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -852995.4399999999
$ws.Range("H100").Value = 1988.9565
$ws.Range("J100").Value = 4105.7144
$ws.Range("L100").Value = 4105.7144
$ws.Range("N100").Value = -5187.7144
$ws.Range("H108").Value = 46300
$ws.Range("J108").Value = 46300
$ws.Range("L108").Value = 46300
$ws.Range("N108").Value = -53980
$ws.Range("H111").Value = 368500
$ws.Range("J111").Value = 368500
$ws.Range("L111").Value = 368500
$ws.Range("N111").Value = -376680
$ws.Range("H114").Value = 12398
$ws.Range("J114").Value = 12398
$ws.Range("L114").Value = 12398
$ws.Range("N114").Value = -21076
$ws.Range("H132").Value = 98049.46000000001
$ws.Range("I132").Value = 98049.46000000001
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 294148.38
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -291618.38
$ws.Range("N132").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 21693.143
$ws.Range("I69").Value = 16001
$ws.Range("J69").Value = 22641.834
$ws.Range("K69").Value = 16001
$ws.Range("L69").Value = 22641.834
$ws.Range("M69").Value = -15252
$ws.Range("N69").Value = -24139.834
$ws.Range("H72").Value = 21693.143
$ws.Range("I72").Value = 16001
$ws.Range("J72").Value = 22641.834
$ws.Range("K72").Value = 48003
$ws.Range("L72").Value = 67925.50199999999
$ws.Range("M72").Value = -44259
$ws.Range("N72").Value = -75413.50199999999
$ws.Range("H81").Value = 3019.4
$ws.Range("I81").Value = 3019.4
$ws.Range("K81").Value = 6038.8
$ws.Range("M81").Value = -4977.8
$ws.Range("H84").Value = 3019.4
$ws.Range("I84").Value = 3019.4
$ws.Range("K84").Value = 30194
$ws.Range("M84").Value = -24890
$ws.Range("H86").Value = 80000
$ws.Range("J86").Value = 80000
$ws.Range("L86").Value = 80000
$ws.Range("N86").Value = -82246
$ws.Range("H89").Value = 80000
$ws.Range("J89").Value = 80000
$ws.Range("L89").Value = 400000
$ws.Range("N89").Value = -411232
$ws.Range("H107").Value = 1355
$ws.Range("I107").Value = 846.2143
$ws.Range("K107").Value = 2538.6429
$ws.Range("M107").Value = -618.6428999999998
$ws.Range("H113").Value = 1450
$ws.Range("I113").Value = 745.5
$ws.Range("J113").Value = 1802.25
$ws.Range("K113").Value = 2236.5
$ws.Range("L113").Value = 5406.75
$ws.Range("M113").Value = -66.5
$ws.Range("N113").Value = -9746.75
$ws.Range("H126").Value = 38100.734
$ws.Range("I126").Value = 44640.36
$ws.Range("K126").Value = 133921.08
$ws.Range("M126").Value = -131451.08
$ws.Range("H127").Value = 99996.5
$ws.Range("J127").Value = 99996.5
$ws.Range("L127").Value = 99996.5
$ws.Range("N127").Value = -109916.5
$ws.Range("H132").Value = 31795.383
$ws.Range("I132").Value = 32667.969
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 98003.90700000001
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -95473.90700000001
$ws.Range("N132").Value = -14060
